$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.306.40'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '2.490.20'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '321.67'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Value = '109.05'
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E10").Value = '  +3.48%  '
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '18.52'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").Value = '2.879.46'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '2.494.78'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '0.845'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '47.236.42'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").Value = '13.38'
$ws.Range("E19").Value = '  +4.46%  '
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +14.47%  '
$ws.Range("D23").Value = '70.62'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = '246.92'
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  +3.78%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("D31").Value = '34.55'
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").Value = '49.87'
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").Value = '20.45'
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").Value = '0.0784'
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '4.76'
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  -2.55%  '
$ws.Range("D40").Value = '22.71'
$ws.Range("E40").Value = '  +6.37%  '
$ws.Range("D41").Value = '0.111'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("E42").Value = '  -2.39%  '
$ws.Range("D43").Value = '119.23'
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = '1.993.91'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").Value = '5.18'
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("D51").Value = '56.77'
$ws.Range("E51").Value = '  +2.89%  '
